$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.759.07'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '1.619.88'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.83'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5087'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2566'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06354'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.27'
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07766'
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.242'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').Value = '1.625.22'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').Value = '1.843.08'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5533'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.50'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('D17').Value = '0.0₅7510'
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('D18').Value = '25.759.08'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.59'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.334'
$ws.Range('E21').Value = '  -3.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.748'
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.965'
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.831'
$ws.Range('E25').Value = '  -3.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.63'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('E28').Value = '  -2.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.40'
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.234'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04855'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.300'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.172'
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.549'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8923'
$ws.Range('E36').Value = '  -2.74%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.125.31'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.532'
$ws.Range('E38').Value = '  -2.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5469'
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01558'
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.565'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7908'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.06'
$ws.Range('E44').Value = '  -2.60%  '
$ws.Range('D45').Value = '1.765.83'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -5.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4420'
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.59'
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05064'
$ws.Range('E49').Value = '  -3.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.522'
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  -0.74%  '

# Reset number format to default (General/Normal style) to avoid leaving the
# "Text" number format applied, while keeping the values as text.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
